# Update CircadiPy cosinor analysis results (sawtooth_05) to the
# republished figures values - re-run of CircaDB / CircadiPy simulations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("E2").Value = 22.90000000000014
$ws.Range("G2").Value = [double]"8.132058471055359e-09"
$ws.Range("H2").Value = [double]"1.055865523118668e-07"
$ws.Range("I2").Value = [double]"3.772537837676282e-13"
$ws.Range("K2").Value = 5.833314818029338
$ws.Range("L2").Value = "[3.6218009721600897, 8.044828663898587]"
$ws.Range("M2").Value = [double]"3.344744516908804e-07"
$ws.Range("N2").Value = [double]"3.344744516908804e-07"
$ws.Range("O2").Value = -1.408842351159387
$ws.Range("P2").Value = "[-1.8365266363327724, -0.9811580659860013]"
$ws.Range("Q2").Value = [double]"2.604640947367898e-10"
$ws.Range("R2").Value = [double]"5.209281894735796e-10"
$ws.Range("S2").Value = 11.39446199346395
$ws.Range("T2").Value = "[10.09099914434941, 12.697924842578495]"
$ws.Range("W2").Value = 5.134734734734767
$ws.Range("X2").Value = 3.575975975975997
$ws.Range("Y2").Value = 6.693493493493537

# ---- Row 3 ----
$ws.Range("E3").Value = 23.41000000000022
$ws.Range("G3").Value = [double]"5.597826313596954e-09"
$ws.Range("H3").Value = [double]"1.055865523118668e-07"
$ws.Range("K3").Value = 6.030334212603734
$ws.Range("L3").Value = "[3.7423579429682565, 8.31831048223921]"
$ws.Range("M3").Value = [double]"3.29720012182122e-07"
$ws.Range("N3").Value = [double]"3.344744516908804e-07"
$ws.Range("O3").Value = 0.1823947686768852
$ws.Range("P3").Value = "[-0.25786846606042246, 0.6226580034141929]"
$ws.Range("Q3").Value = 0.4160221521485334
$ws.Range("R3").Value = 0.4160221521485334
$ws.Range("S3").Value = 11.71051953952427
$ws.Range("T3").Value = "[10.379081356043187, 13.04195772300536]"
$ws.Range("W3").Value = 22.73043043043064
$ws.Range("X3").Value = 21.09009009009029
$ws.Range("Y3").Value = 24.370770770771
